$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: mark D-column cells whose new values look like plain numbers as Text, ---
# --- so Excel stores the exact original string instead of converting to a number. ---
$textFormatCells = @(
    "D5", "D7", "D9", "D10", "D11", "D12", "D13", "D14", "D16", "D17", "D21", "D22", "D24", "D25", "D28", "D29", "D34", "D38", "D39", "D41", "D47", "D48", "D51"
)
foreach ($addr in $textFormatCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# --- Step 2: write the new cell values row by row ---
# Row 2
$ws.Range("D2").Value = "36.483.76"
$ws.Range("E2").Value = "  -1.17%  "
# Row 3
$ws.Range("D3").Value = "1.962.62"
$ws.Range("E3").Value = "  -3.75%  "
# Row 4
$ws.Range("E4").Value = "  +0.12%  "
# Row 5
$ws.Range("D5").Value = "244.36"
$ws.Range("E5").Value = "  -1.88%  "
# Row 6
$ws.Range("E6").Value = "  -2.95%  "
# Row 7
$ws.Range("D7").Value = "58.88"
$ws.Range("E7").Value = "  -6.64%  "
# Row 8
$ws.Range("E8").Value = "  +0.05%  "
# Row 9
$ws.Range("D9").Value = "0.376"
$ws.Range("E9").Value = "  -3.18%  "
# Row 10
$ws.Range("D10").Value = "55.84"
$ws.Range("E10").Value = "  -4.07%  "
# Row 11
$ws.Range("D11").Value = "0.0846"
$ws.Range("E11").Value = "  +6.23%  "
# Row 12
$ws.Range("D12").Value = "0.104"
$ws.Range("E12").Value = "  -0.07%  "
# Row 13
$ws.Range("B13").Value = "Polygon"
$ws.Range("C13").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D13").Value = "0.841"
$ws.Range("E13").Value = "  -7.17%  "
# Row 14
$ws.Range("B14").Value = "Avalanche"
$ws.Range("C14").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D14").Value = "22.02"
$ws.Range("E14").Value = "  -4.91%  "
# Row 15
$ws.Range("D15").Value = "2.253.87"
$ws.Range("E15").Value = "  -3.57%  "
# Row 16
$ws.Range("D16").Value = "13.62"
$ws.Range("E16").Value = "  -5.31%  "
# Row 17
$ws.Range("D17").Value = "5.37"
$ws.Range("E17").Value = "  -3.34%  "
# Row 18
$ws.Range("D18").Value = "1.981.62"
$ws.Range("E18").Value = "  -2.77%  "
# Row 19
$ws.Range("D19").Value = "36.401.96"
# Row 20
$ws.Range("D20").Value = "0.0₃0885"
$ws.Range("E20").Value = "  +0.05%  "
# Row 21
$ws.Range("D21").Value = "70.38"
$ws.Range("E21").Value = "  -2.63%  "
# Row 22
$ws.Range("D22").Value = "231.54"
$ws.Range("E22").Value = "  -2.22%  "
# Row 23
$ws.Range("E23").Value = "  -5.55%  "
# Row 24
$ws.Range("D24").Value = "1.00"
$ws.Range("E24").Value = "  +0.04%  "
# Row 25
$ws.Range("D25").Value = "2.53"
$ws.Range("E25").Value = "  +0.29%  "
# Row 26
$ws.Range("E26").Value = "  -2.40%  "
# Row 27
$ws.Range("E27").Value = "  -2.20%  "
# Row 28
$ws.Range("D28").Value = "165.03"
$ws.Range("E28").Value = "  +3.38%  "
# Row 29
$ws.Range("D29").Value = "19.73"
$ws.Range("E29").Value = "  -2.77%  "
# Row 30
$ws.Range("E30").Value = "  -14.22%  "
# Row 31
$ws.Range("E31").Value = "  -2.07%  "
# Row 32
$ws.Range("E32").Value = "  -0.92%  "
# Row 33
$ws.Range("E33").Value = "  -5.79%  "
# Row 34
$ws.Range("D34").Value = "0.0642"
$ws.Range("E34").Value = "  +3.45%  "
# Row 35
$ws.Range("E35").Value = "  -3.31%  "
# Row 36
$ws.Range("E36").Value = "  -2.28%  "
# Row 37
$ws.Range("E37").Value = "  -0.01%  "
# Row 38
$ws.Range("D38").Value = "1.81"
$ws.Range("E38").Value = "  -1.86%  "
# Row 39
$ws.Range("D39").Value = "2.17"
$ws.Range("E39").Value = "  -8.69%  "
# Row 40
$ws.Range("E40").Value = "  -6.92%  "
# Row 41
$ws.Range("D41").Value = "0.0985"
$ws.Range("E41").Value = "  -1.31%  "
# Row 42
$ws.Range("E42").Value = "  -4.32%  "
# Row 43
$ws.Range("E43").Value = "  -3.25%  "
# Row 44
$ws.Range("E44").Value = "  -1.72%  "
# Row 45
$ws.Range("E45").Value = "  -7.14%  "
# Row 46
$ws.Range("E46").Value = "  -7.41%  "
# Row 47
$ws.Range("D47").Value = "7.45"
$ws.Range("E47").Value = "  -3.39%  "
# Row 48
$ws.Range("D48").Value = "89.35"
$ws.Range("E48").Value = "  -4.87%  "
# Row 49
$ws.Range("D49").Value = "1.350.18"
$ws.Range("E49").Value = "  -1.49%  "
# Row 50
$ws.Range("E50").Value = "  -3.22%  "
# Row 51
$ws.Range("D51").Value = "45.38"
$ws.Range("E51").Value = "  -0.36%  "

# --- Step 3: restore Normal style on the text-formatted D cells so no stray
# --- cell-level style index is left behind (matches original formatting). ---
foreach ($addr in $textFormatCells) {
    $ws.Range($addr).Style = "Normal"
}
